$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-13 with new TPM-based values.
# Columns B, C, E, F remain unchanged (Ligand symbol, Receptor symbol, Ligand-expressing cells, Ligand detection rate).

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 4.529578333333333
$ws.Range("H2").Value = 13.588735
$ws.Range("I2").Value = 0.2308013058217703
$ws.Range("J2").Value = 0.2308013058217703
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.5447316666666667
$ws.Range("N2").Value = 1.634195
$ws.Range("O2").Value = 0.1484165462704666
$ws.Range("P2").Value = 0.1484165462704666
$ws.Range("Q2").Value = 2.467404754813889
$ws.Range("R2").Value = 22.206642793325
$ws.Range("S2").Value = 0.03425473268478089
$ws.Range("T2").Value = 0.03425473268478089

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 4.529578333333333
$ws.Range("H3").Value = 13.588735
$ws.Range("I3").Value = 0.2308013058217703
$ws.Range("J3").Value = 0.2308013058217703
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.371854333333333
$ws.Range("N3").Value = 7.115563
$ws.Range("O3").Value = 0.6462308875194944
$ws.Range("P3").Value = 0.6462308875194943
$ws.Range("Q3").Value = 10.74349999808944
$ws.Range("R3").Value = 96.691499982805
$ws.Range("S3").Value = 0.1491509327018609
$ws.Range("T3").Value = 0.1491509327018609

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 4.529578333333333
$ws.Range("H4").Value = 13.588735
$ws.Range("I4").Value = 0.2308013058217703
$ws.Range("J4").Value = 0.2308013058217703
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.7537033333333333
$ws.Range("N4").Value = 2.26111
$ws.Range("O4").Value = 0.205352566210039
$ws.Range("P4").Value = 0.205352566210039
$ws.Range("Q4").Value = 3.413958288427777
$ws.Range("R4").Value = 30.72562459585
$ws.Range("S4").Value = 0.04739564043512856
$ws.Range("T4").Value = 0.04739564043512856

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "ECs"
$ws.Range("G5").Value = 6.782643666666666
$ws.Range("H5").Value = 20.347931
$ws.Range("I5").Value = 0.3456045794970084
$ws.Range("J5").Value = 0.3456045794970085
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.5447316666666667
$ws.Range("N5").Value = 1.634195
$ws.Range("O5").Value = 0.1484165462704666
$ws.Range("P5").Value = 0.1484165462704666
$ws.Range("Q5").Value = 3.694720788949445
$ws.Range("R5").Value = 33.252487100545
$ws.Range("S5").Value = 0.0512934380642029
$ws.Range("T5").Value = 0.05129343806420291

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("G6").Value = 6.782643666666666
$ws.Range("H6").Value = 20.347931
$ws.Range("I6").Value = 0.3456045794970084
$ws.Range("J6").Value = 0.3456045794970085
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.371854333333333
$ws.Range("N6").Value = 7.115563
$ws.Range("O6").Value = 0.6462308875194944
$ws.Range("P6").Value = 0.6462308875194943
$ws.Range("Q6").Value = 16.08744277223922
$ws.Range("R6").Value = 144.786984950153
$ws.Range("S6").Value = 0.2233403541391534
$ws.Range("T6").Value = 0.2233403541391534

# Row 7: FAPs -> MuSCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("G7").Value = 6.782643666666666
$ws.Range("H7").Value = 20.347931
$ws.Range("I7").Value = 0.3456045794970084
$ws.Range("J7").Value = 0.3456045794970085
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7537033333333333
$ws.Range("N7").Value = 2.26111
$ws.Range("O7").Value = 0.205352566210039
$ws.Range("P7").Value = 0.205352566210039
$ws.Range("Q7").Value = 5.112101140378888
$ws.Range("R7").Value = 46.00891026340999
$ws.Range("S7").Value = 0.07097078729365211
$ws.Range("T7").Value = 0.07097078729365212

# Row 8: MuSCs -> ECs
$ws.Range("A8").Value = "MuSCs"
$ws.Range("D8").Value = "ECs"
$ws.Range("G8").Value = 2.766332333333333
$ws.Range("H8").Value = 8.298997
$ws.Range("I8").Value = 0.1409564131327128
$ws.Range("J8").Value = 0.1409564131327128
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.5447316666666667
$ws.Range("N8").Value = 1.634195
$ws.Range("O8").Value = 0.1484165462704666
$ws.Range("P8").Value = 0.1484165462704666
$ws.Range("Q8").Value = 1.506908822490556
$ws.Range("R8").Value = 13.562179402415
$ws.Range("S8").Value = 0.02092026401183028
$ws.Range("T8").Value = 0.02092026401183028

# Row 9: MuSCs -> FAPs
$ws.Range("A9").Value = "MuSCs"
$ws.Range("D9").Value = "FAPs"
$ws.Range("G9").Value = 2.766332333333333
$ws.Range("H9").Value = 8.298997
$ws.Range("I9").Value = 0.1409564131327128
$ws.Range("J9").Value = 0.1409564131327128
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.371854333333333
$ws.Range("N9").Value = 7.115563
$ws.Range("O9").Value = 0.6462308875194944
$ws.Range("P9").Value = 0.6462308875194943
$ws.Range("Q9").Value = 6.561337332256778
$ws.Range("R9").Value = 59.052035990311
$ws.Range("S9").Value = 0.09109038796031753
$ws.Range("T9").Value = 0.09109038796031753

# Row 10: MuSCs -> MuSCs
$ws.Range("A10").Value = "MuSCs"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("G10").Value = 2.766332333333333
$ws.Range("H10").Value = 8.298997
$ws.Range("I10").Value = 0.1409564131327128
$ws.Range("J10").Value = 0.1409564131327128
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.7537033333333333
$ws.Range("N10").Value = 2.26111
$ws.Range("O10").Value = 0.205352566210039
$ws.Range("P10").Value = 0.205352566210039
$ws.Range("Q10").Value = 2.084993900741111
$ws.Range("R10").Value = 18.76494510667
$ws.Range("S10").Value = 0.02894576116056502
$ws.Range("T10").Value = 0.02894576116056503

# Row 11: Resolving-Mac -> ECs
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("D11").Value = "ECs"
$ws.Range("G11").Value = 5.546890666666666
$ws.Range("H11").Value = 16.640672
$ws.Range("I11").Value = 0.2826377015485084
$ws.Range("J11").Value = 0.2826377015485084
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.5447316666666667
$ws.Range("N11").Value = 1.634195
$ws.Range("O11").Value = 0.1484165462704666
$ws.Range("P11").Value = 0.1484165462704666
$ws.Range("Q11").Value = 3.021566997671111
$ws.Range("R11").Value = 27.19410297904
$ws.Range("S11").Value = 0.04194811150965253
$ws.Range("T11").Value = 0.04194811150965253

# Row 12: Resolving-Mac -> FAPs
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("D12").Value = "FAPs"
$ws.Range("G12").Value = 5.546890666666666
$ws.Range("H12").Value = 16.640672
$ws.Range("I12").Value = 0.2826377015485084
$ws.Range("J12").Value = 0.2826377015485084
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.371854333333333
$ws.Range("N12").Value = 7.115563
$ws.Range("O12").Value = 0.6462308875194944
$ws.Range("P12").Value = 0.6462308875194943
$ws.Range("Q12").Value = 13.15641666425955
$ws.Range("R12").Value = 118.407749978336
$ws.Range("S12").Value = 0.1826492127181626
$ws.Range("T12").Value = 0.1826492127181626

# Row 13: Resolving-Mac -> MuSCs
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("G13").Value = 5.546890666666666
$ws.Range("H13").Value = 16.640672
$ws.Range("I13").Value = 0.2826377015485084
$ws.Range("J13").Value = 0.2826377015485084
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.7537033333333333
$ws.Range("N13").Value = 2.26111
$ws.Range("O13").Value = 0.205352566210039
$ws.Range("P13").Value = 0.205352566210039
$ws.Range("Q13").Value = 4.180709985102221
$ws.Range("R13").Value = 37.62638986592
$ws.Range("S13").Value = 0.05804037732069332
$ws.Range("T13").Value = 0.05804037732069332

# Remove the trailing rows (14:17) that corresponded to the removed "Resolving-Mac" sending cluster block.
$ws.Rows("14:17").Delete()
